$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full symmetric 24x24 distance matrix (rows/cols 1-20 original, 21-24 new
# locations added for Problem 2 calibration and demand prediction).
$rows = @(
@(0,347,370,654,1246,1148,941,1444,449,1462,1564,947,1848,1470,533,2104,2715,1895,1834,2472,5541,7764,8762,6345),
@(347,0,398,449,1064,859,681,1101,785,1539,1470,850,1896,1342,868,1850,2390,2241,1487,2422,5835,8073,9105,6666),
@(370,398,0,366,1461,1241,664,1297,750,1152,1846,578,1521,1101,666,1783,2481,2038,1712,2784,5849,8051,8958,6613),
@(654,449,366,0,1422,1093,299,958,1087,1223,1874,431,1538,896,1033,1449,2115,2401,1378,2836,6191,8404,9324,6973),
@(1246,1064,1461,1422,0,483,1497,1330,1454,2602,513,1854,2948,2271,1721,2458,2574,2893,1473,1460,5763,8041,9389,6747),
@(1148,859,1241,1093,483,0,1094,847,1485,2315,994,1510,2628,1870,1677,1981,2114,2982,1023,1922,6152,8427,9677,7095),
@(941,681,664,299,1497,1094,0,729,1381,1317,1984,479,1575,777,1330,1171,1818,2692,1145,2948,6483,8700,9618,7271),
@(1444,1101,1297,958,1330,847,729,0,1886,2023,1839,1199,2234,1327,1934,1162,1333,3320,420,2744,6867,9127,10206,7744),
@(449,785,750,1087,1454,1485,1381,1886,0,1625,1642,1308,2023,1826,337,2530,3164,1498,2267,2451,5104,7319,8320,5897),
@(1462,1539,1152,1223,2602,2315,1317,2023,1625,0,2999,839,398,853,1320,1764,2754,2143,2418,3930,6293,8343,8865,6857),
@(1564,1470,1846,1874,513,994,1984,1839,1642,2999,0,2303,3364,2749,1956,2971,3051,2941,1952,965,5406,7678,9139,6437),
@(947,850,578,431,1854,1510,479,1199,1308,839,2303,0,1119,523,1142,1291,2132,2406,1607,3263,6369,8540,9321,7083),
@(1848,1896,1521,1538,2948,2628,1575,2234,2023,398,3364,1119,0,939,1718,1753,2778,2445,2602,4306,6609,8610,9016,7121),
@(1470,1342,1101,896,2271,1870,777,1327,1826,853,2749,523,939,0,1634,924,1901,2805,1672,3715,6849,8990,9661,7521),
@(533,868,666,1033,1721,1677,1330,1934,337,1320,1956,1142,1718,1634,0,2419,3148,1385,2339,2786,5226,7405,8294,5958),
@(2104,1850,1783,1449,2458,1981,1171,1162,2530,1764,2971,1291,1753,924,2419,0,1029,3690,1283,3899,7633,9825,10580,8373),
@(2715,2390,2481,2115,2574,2114,1818,1333,3164,2754,3051,2132,2778,1901,3148,1029,0,4506,1101,3850,8199,10454,11435,9056),
@(1895,2241,2038,2401,2893,2982,2692,3320,1498,2143,2941,2406,2445,2805,1385,3690,4506,0,3723,3503,4164,6207,6928,4726),
@(1834,1487,1712,1378,1473,1023,1145,420,2267,2418,1952,1607,2602,1672,2339,1283,1101,3723,0,2780,7161,9432,10583,8075),
@(2472,2422,2784,2836,1460,1922,2948,2744,2451,3930,965,3263,4306,3715,2786,3899,3850,3503,2780,0,5072,7286,8981,6194),
@(5541,5835,5849,6191,5763,6152,6483,6867,5104,6293,5406,6369,6609,6849,5226,7633,8199,4164,7161,5072,0,2277,3975,1188),
@(7764,8073,8051,8404,8041,8427,8700,9127,7319,8343,7678,8540,8610,8990,7405,9825,10454,6207,9432,7286,2277,0,2216,1490),
@(8762,9105,8958,9324,9389,9677,9618,10206,8320,8865,9139,9321,9016,9661,8294,10580,11435,6928,10583,8981,3975,2216,0,2802),
@(6345,6666,6613,6973,6747,7095,7271,7744,5897,6857,6437,7083,7121,7521,5958,8373,9056,4726,8075,6194,1188,1490,2802,0)
)

$data = New-Object 'object[,]' $rows.Count, $rows[0].Count
for ($i = 0; $i -lt $rows.Count; $i++) {
    for ($j = 0; $j -lt $rows[$i].Count; $j++) {
        $data[$i, $j] = $rows[$i][$j]
    }
}

$ws.Range("A1:X24").Value = $data
